$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New "XP" doc-number columns (F, G) + SQL-insert helper columns
# (J, K, O, P, Q) for rows 2-8. This mirrors a fill-down of a doc
# reference ("01".."07"), a constant prefix ("XP"), and formulas that
# stitch together a SQL INSERT INTO supplier_debts(...) VALUES (...)
# statement out of the existing row data.
# ---------------------------------------------------------------------

# F2:F8 - two-digit running reference, kept as text ("01".."07")
for ($r = 2; $r -le 8; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.NumberFormat = "@"
    $cell.Value = "0$($r - 1)"
}

# G2:G8 - constant "XP" prefix, styled like the other text columns (A:D)
$ws.Range("A2").Copy()
$ws.Range("G2:G8").PasteSpecial(-4122)
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 7).Value = "XP"
}

# J2:J8 - literal start of the INSERT statement
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 10).Value = "INSERT INTO supplier_debts (reference,total,doc_no,type_debt) VALUES ("
}

# Row 2 formulas (entered individually, so they stay standalone formulas)
$ws.Range("K2").Formula = "=""'""&A2&""',"""
$ws.Range("O2").Formula = "=""'""&E2&""',"""
$ws.Range("P2").Formula = "=""'""&G2&""2006-000""&F2&""',"""
$ws.Range("Q2").Formula = "=""'""&G2&""');"""

# Rows 3-8 formulas (filled down as one range => shared formulas)
$ws.Range("K3:K8").Formula = "=""'""&A3&""',"""
$ws.Range("O3:O8").Formula = "=""'""&E3&""',"""
$ws.Range("P3:P8").Formula = "=""'""&G3&""2006-000""&F3&""',"""
$ws.Range("Q3:Q8").Formula = "=""'""&G3&""');"""

# ---------------------------------------------------------------------
# Approximate the column widths of the brand-new columns so they are
# not left at the bare default. (Existing columns A:E keep their
# original bestFit widths untouched.)
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 8.3349
$ws.Columns.Item(7).ColumnWidth = 17.6701
$ws.Columns.Item(16).ColumnWidth = 13.5026
$ws.Columns.Item(17).ColumnWidth = 16.6699

# ---------------------------------------------------------------------
# Selection / view: user ends up with J3:Q8 selected (active cell J3)
# after scrolling the sheet so column C is the left-most visible one.
# ---------------------------------------------------------------------
$ws.Range("J3:Q8").Select()
